$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting existing D:K data to F:M
$ws.Range("D:E").Insert()

# Copy cell formatting (number format/font/style) from the shifted columns (F:G)
# into the newly inserted D:E columns, per contiguous data block, so that
# each row keeps its original per-row style (date header rows vs numeric rows).
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new D:E columns with the new quarter figures
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 2627000
$ws.Range("E8").Value = 2524000
$ws.Range("D9").Value = 1484000
$ws.Range("E9").Value = 1340000
$ws.Range("D10").Value = 1143000
$ws.Range("E10").Value = 1184000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 19000
$ws.Range("E14").Value = 18000
$ws.Range("D15").Value = 163000
$ws.Range("E15").Value = 163000
$ws.Range("D17").Value = 2453000
$ws.Range("E17").Value = 2346000
$ws.Range("D18").Value = 174000
$ws.Range("E18").Value = 178000
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 337000
$ws.Range("E21").Value = 341000
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 174000
$ws.Range("E23").Value = 178000
$ws.Range("D24").Value = 55000
$ws.Range("E24").Value = 50000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 119000
$ws.Range("E26").Value = 128000
$ws.Range("D27").Value = 95000
$ws.Range("E27").Value = 101000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 95000
$ws.Range("E33").Value = 101000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 95000
$ws.Range("E35").Value = 101000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 1618000
$ws.Range("E41").Value = 1886000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 1853000
$ws.Range("E43").Value = 1648000
$ws.Range("D44").Value = 400000
$ws.Range("E44").Value = 388000
$ws.Range("D45").Value = 558000
$ws.Range("E45").Value = 547000
$ws.Range("D46").Value = 4429000
$ws.Range("E46").Value = 4469000
$ws.Range("D47").Value = 345000
$ws.Range("E47").Value = 390000
$ws.Range("D48").Value = 2517000
$ws.Range("E48").Value = 2512000
$ws.Range("D49").Value = 7796000
$ws.Range("E49").Value = 7760000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1140000
$ws.Range("E52").Value = 1157000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 16227000
$ws.Range("E54").Value = 16288000
$ws.Range("D57").Value = 625000
$ws.Range("E57").Value = 537000
$ws.Range("D58").Value = 744000
$ws.Range("E58").Value = 671000
$ws.Range("D59").Value = 2343000
$ws.Range("E59").Value = 2337000
$ws.Range("D60").Value = 3712000
$ws.Range("E60").Value = 3545000
$ws.Range("D61").Value = 936000
$ws.Range("E61").Value = 1186000
$ws.Range("D62").Value = 1145000
$ws.Range("E62").Value = 1127000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 6963000
$ws.Range("E66").Value = 7027000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -1937000
$ws.Range("E72").Value = -2032000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 9264000
$ws.Range("E76").Value = 9261000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 95000
$ws.Range("E81").Value = 101000
$ws.Range("D83").Value = 163000
$ws.Range("E83").Value = 163000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 245000
$ws.Range("E89").Value = 113000
$ws.Range("D91").Value = -131000
$ws.Range("E91").Value = -133000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -288000
$ws.Range("E94").Value = -121000
$ws.Range("D96").Value = -58000
$ws.Range("E96").Value = -23000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -209000
$ws.Range("E100").Value = -124000
$ws.Range("D101").Value = -16000
$ws.Range("E101").Value = -16000
$ws.Range("D102").Value = -268000
$ws.Range("E102").Value = -148000
